$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.056.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = "'1.667.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'216.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("D6").Value = "'0.5101"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D9").Value = "'0.06392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("D10").Value = "'21.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").Value = "'0.07448"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = "'1.681.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").Value = "'4.511"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = "'0.5808"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = "'0.000008494"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").Value = "'64.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = "'25.887.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.21%  '
$ws.Range("D18").Value = "'4.928"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").Value = "'189.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").Value = "'6.191"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = "'145.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = "'7.617"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = "'0.1218"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.88%  '
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = "'0.06687"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.22%  '
$ws.Range("D29").Value = "'1.331"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").Value = "'1.313"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = "'3.556"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("D32").Value = "'3.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = "'1.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = "'1.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").Value = "'0.6167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.33%  '
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = "'2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").Value = "'6.325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.23%  '
$ws.Range("D39").Value = "'1.097.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").Value = "'0.8695"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").Value = "'1.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("D43").Value = "'101.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").Value = "'1.816.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").Value = "'0.00000000115"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.24%  '
$ws.Range("D46").Value = "'56.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = "'1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = "'8.117"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  +0.28%  '
$ws.Range("D50").Value = "'0.4277"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").Value = "'5.990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
